$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# are pre-formatted as Text ("@") so the literal string content is preserved.

$ws.Range("D2").Value = "31.071.32"
$ws.Range("E2").Value = "  +1.63%  "
$ws.Range("D3").Value = "1.959.27"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  +0.65%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.02"
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("E7").Value = "  +1.85%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2950"
$ws.Range("E8").Value = "  +1.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07020"
$ws.Range("E9").Value = "  +4.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.42"
$ws.Range("E10").Value = "  +1.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "107.56"
$ws.Range("E11").Value = "  -1.38%  "
$ws.Range("D12").Value = "1.955.45"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07796"
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.493"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7005"
$ws.Range("E15").Value = "  +1.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "282.66"
$ws.Range("E16").Value = "  -2.89%  "
$ws.Range("D17").Value = "31.093.51"
$ws.Range("E17").Value = "  +1.62%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007782"
$ws.Range("E18").Value = "  +1.61%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.25"
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("D20").Value = "2.217.62"
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.519"
$ws.Range("E22").Value = "  -1.80%  "
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.508"
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.854"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.41"
$ws.Range("E26").Value = "  -1.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.98"
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.192"
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1049"
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.393"
$ws.Range("E30").Value = "  -3.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.579"
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.621"
$ws.Range("E32").Value = "  -2.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.427"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04926"
$ws.Range("E34").Value = "  -3.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7553"
$ws.Range("E35").Value = "  -1.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.170"
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.735"
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02007"
$ws.Range("E38").Value = "  -1.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.705"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.530"
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "78.01"
$ws.Range("E41").Value = "  +11.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.120"
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9037"
$ws.Range("E43").Value = "  +2.65%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4460"
$ws.Range("E44").Value = "  +0.46%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "109.28"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.142"
$ws.Range("E46").Value = "  +9.08%  "
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("D48").Value = "1.026.51"
$ws.Range("E48").Value = "  +10.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.330"
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1256"
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.94"
$ws.Range("E51").Value = "  +0.28%  "
